$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '25.700.83'
Set-TextCell $ws.Range('E2') '  -3.44%  '
Set-TextCell $ws.Range('D3') '1.744.21'
Set-TextCell $ws.Range('E3') '  -5.73%  '
Set-TextCell $ws.Range('D4') '1.001'
Set-TextCell $ws.Range('E4') '  -0.05%  '
Set-TextCell $ws.Range('D5') '236.26'
Set-TextCell $ws.Range('E5') '  -10.11%  '
Set-TextCell $ws.Range('D6') '1.001'
Set-TextCell $ws.Range('E6') '  +0.00%  '
Set-TextCell $ws.Range('D7') '0.4935'
Set-TextCell $ws.Range('E7') '  -7.51%  '
Set-TextCell $ws.Range('D8') '41.37'
Set-TextCell $ws.Range('E8') '  -8.50%  '
Set-TextCell $ws.Range('D9') '0.2562'
Set-TextCell $ws.Range('E9') '  -18.88%  '
Set-TextCell $ws.Range('D10') '0.06030'
Set-TextCell $ws.Range('E10') '  -13.11%  '
Set-TextCell $ws.Range('D11') '1.745.04'
Set-TextCell $ws.Range('E11') '  -5.71%  '
Set-TextCell $ws.Range('D12') '0.06831'
Set-TextCell $ws.Range('E12') '  -12.84%  '
Set-TextCell $ws.Range('D13') '14.83'
Set-TextCell $ws.Range('E13') '  -21.31%  '
Set-TextCell $ws.Range('D14') '4.448'
Set-TextCell $ws.Range('E14') '  -11.95%  '
Set-TextCell $ws.Range('D15') '76.61'
Set-TextCell $ws.Range('E15') '  -14.64%  '
Set-TextCell $ws.Range('D16') '0.5714'
Set-TextCell $ws.Range('E16') '  -25.81%  '
Set-TextCell $ws.Range('D17') '1.000'
Set-TextCell $ws.Range('E17') '  -0.16%  '
Set-TextCell $ws.Range('E18') '  -0.06%  '
Set-TextCell $ws.Range('D19') '25.730.11'
Set-TextCell $ws.Range('E19') '  -3.41%  '
Set-TextCell $ws.Range('D20') '11.30'
Set-TextCell $ws.Range('E20') '  -19.99%  '
Set-TextCell $ws.Range('D21') '0.000006565'
Set-TextCell $ws.Range('E21') '  -17.58%  '
Set-TextCell $ws.Range('D22') '1.967.10'
Set-TextCell $ws.Range('E22') '  -5.84%  '
Set-TextCell $ws.Range('D23') '4.012'
Set-TextCell $ws.Range('E23') '  -13.75%  '
Set-TextCell $ws.Range('D24') '5.070'
Set-TextCell $ws.Range('E24') '  -15.85%  '
Set-TextCell $ws.Range('D25') '7.940'
Set-TextCell $ws.Range('E25') '  -15.14%  '
Set-TextCell $ws.Range('D26') '137.32'
Set-TextCell $ws.Range('E26') '  -3.14%  '
Set-TextCell $ws.Range('D27') '1.474'
Set-TextCell $ws.Range('E27') '  -12.68%  '
Set-TextCell $ws.Range('D28') '1.823'
Set-TextCell $ws.Range('E28') '  -17.56%  '
Set-TextCell $ws.Range('D29') '14.67'
Set-TextCell $ws.Range('E29') '  -14.36%  '
Set-TextCell $ws.Range('D30') '101.92'
Set-TextCell $ws.Range('E30') '  -8.79%  '
Set-TextCell $ws.Range('D31') '3.763'
Set-TextCell $ws.Range('E31') '  -12.66%  '
Set-TextCell $ws.Range('D32') '0.07987'
Set-TextCell $ws.Range('E32') '  -9.07%  '
Set-TextCell $ws.Range('D33') '3.415'
Set-TextCell $ws.Range('E33') '  -17.08%  '
Set-TextCell $ws.Range('D34') '0.04387'
Set-TextCell $ws.Range('E34') '  -9.61%  '
Set-TextCell $ws.Range('D35') '0.9993'
Set-TextCell $ws.Range('E35') '  -0.10%  '
Set-TextCell $ws.Range('D36') '2.610'
Set-TextCell $ws.Range('E36') '  -9.90%  '
Set-TextCell $ws.Range('D37') '0.9832'
Set-TextCell $ws.Range('E37') '  -13.68%  '
Set-TextCell $ws.Range('D38') '0.5992'
Set-TextCell $ws.Range('E38') '  -18.84%  '
Set-TextCell $ws.Range('D39') '2.667'
Set-TextCell $ws.Range('E39') '  -14.37%  '
Set-TextCell $ws.Range('D40') '1.962'
Set-TextCell $ws.Range('E40') '  -16.10%  '
Set-TextCell $ws.Range('D41') '1.001'
Set-TextCell $ws.Range('E41') '  +0.01%  '
Set-TextCell $ws.Range('E42') '  -12.92%  '
Set-TextCell $ws.Range('D43') '101.84'
Set-TextCell $ws.Range('E43') '  -6.20%  '
Set-TextCell $ws.Range('D44') '0.7557'
Set-TextCell $ws.Range('D45') '5.153'
Set-TextCell $ws.Range('E45') '  -12.82%  '
Set-TextCell $ws.Range('D46') '0.3765'
Set-TextCell $ws.Range('E46') '  -22.03%  '
Set-TextCell $ws.Range('D47') '0.05228'
Set-TextCell $ws.Range('E47') '  -10.05%  '
Set-TextCell $ws.Range('D48') '0.1069'
Set-TextCell $ws.Range('E48') '  -14.46%  '
Set-TextCell $ws.Range('D49') '30.14'
Set-TextCell $ws.Range('E49') '  -13.91%  '
Set-TextCell $ws.Range('D50') '52.23'
Set-TextCell $ws.Range('E50') '  -13.60%  '
Set-TextCell $ws.Range('D51') '5.825'
Set-TextCell $ws.Range('E51') '  -24.25%  '
